$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: populate the first product / stock row ---

# A7: "م" (row number) - plain numeric value.
$ws.Range("A7").Value = 1

# C7:G7 ("الاسم" / name) and N7:O7 ("السعر" / price) share one style in
# the template, so switch every one of those cells to Text format
# together (cell by cell, so they all land on the same regenerated
# style) before filling in the two values that actually carry text.
$ws.Cells.Item(7, 3).NumberFormat = "@"   # C7
$ws.Cells.Item(7, 4).NumberFormat = "@"   # D7
$ws.Cells.Item(7, 5).NumberFormat = "@"   # E7
$ws.Cells.Item(7, 6).NumberFormat = "@"   # F7
$ws.Cells.Item(7, 7).NumberFormat = "@"   # G7
$ws.Cells.Item(7, 14).NumberFormat = "@"  # N7
$ws.Cells.Item(7, 15).NumberFormat = "@"  # O7

$ws.Range("C7").Value = "سرنجات 3 سم"
$ws.Range("N7").Value = "2.00"

# H7:K7 ("الرصيد الحالي" / current balance) - own shared style -> Text.
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("H7").Value = "0:0"

# L7 ("حد الطلب" / order limit) keeps its original numeric format
# (numFmtId 165) but now holds a text value, so flip to Text only long
# enough to store the string, then restore the original format.
$fmtL7 = $ws.Range("L7").NumberFormat
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "0"
$ws.Range("L7").NumberFormat = $fmtL7

# P7 ("سعر البيع" / selling price) keeps its original numeric format
# (numFmtId 2) but now holds a text value - same trick as L7.
$fmtP7 = $ws.Range("P7").NumberFormat
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "2.0000"
$ws.Range("P7").NumberFormat = $fmtP7

# Q7 ("عدد التعاملات" / transactions count) - own shared style -> Text.
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "1:0"

# --- Row 8: counter cell next to P8:Q8 ---
$ws.Range("P8").Value = 2

# --- Row 9: footer timestamp refreshed to the new save time ---
$ws.Range("A9").Value = "Friday, 11 July, 2025 2:53 PM"
